$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# G2: Area formula for row 2 (uses 0 instead of D1 since it's the first segment)
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
# H2: total area sum
$ws.Range("H2").Formula = "=SUM(G2:G11)"
# J2/K2: summary cells
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# G3
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# G4:G15 shared formula (fill down)
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Match final selection state from the authored session
$ws.Range("J2:K2").Select() | Out-Null

$wb.Save()
